# Updates the Monte Carlo "Model" sheet's simulated deal outcomes
# (Capital out / Date In / Date Out) for companies 1-10 (rows 16-25)
# with freshly re-rolled simulation values. Dependent formulas
# (Duration, MOIC, IRR, Average Hold, totals, portfolio MOIC/IRR)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

$data = @(
    @{ Row = 16; D = 35.993671079192097;  E = 45712; F = 47115 },
    @{ Row = 17; D = 268.76625466002241;  E = 45334; F = 47533 },
    @{ Row = 18; D = 107.45787495594396;  E = 45762; F = 47648 },
    @{ Row = 19; D = 17.583127519126556;  E = 45236; F = 46337 },
    @{ Row = 20; D = 101.54818637999711;  E = 45860; F = 47917 },
    @{ Row = 21; D = 137.471168681488;    E = 45825; F = 48239 },
    @{ Row = 22; D = 210.13320332399107;  E = 45759; F = 47754 },
    @{ Row = 23; D = 7.4955300174986341;  E = 45512; F = 48605 },
    @{ Row = 24; D = 64.264035171550873;  E = 45008; F = 47205 },
    @{ Row = 25; D = 35.251036294235412;  E = 45946; F = 47911 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D   # Column D: Capital out
    $ws.Cells.Item($r, 5).Value = $item.E   # Column E: Date In
    $ws.Cells.Item($r, 6).Value = $item.F   # Column F: Date Out
}

$excel.CalculateFullRebuild()
